$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 1, shifting all existing rows down
$ws.Rows.Item(1).Insert()

# Set the header value for the newly inserted row
$ws.Range("A1").Value = "Symbol"

# Update selection to A2 (as reflected in the saved file)
$ws.Range("A2").Select()
